$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns I (I0) and J (IF)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold font, border, center/top alignment)
# from the existing H1 header cell onto the two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for the new I0 / IF columns, keyed by row number.
$data = @{
    2  = @(8, 8)
    3  = @(7, 7)
    4  = @(7, 8)
    5  = @(7, 8)
    6  = @(5, 5)
    7  = @(7, 8)
    8  = @(1, 1)
    9  = @(8, 8)
    10 = @(8, 8)
    11 = @(7, 7)
    12 = @(7, 7)
    13 = @(8, 8)
    14 = @(9, 9)
    15 = @(9, 9)
    16 = @(5, 5)
    17 = @(7, 7)
    18 = @(5, 5)
    19 = @(5, 5)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
